# Actualización desde MV -datos-
# Updates the revised "01-01-2021" quarter figures (row 102) and appends
# the new "01-04-2021" quarter (row 103).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the label for the new quarter as a genuine text value (not an
# auto-converted date serial). Writing it through a formula and then
# collapsing the formula to its literal value keeps it stored as a
# shared string without introducing any new cell style.
$ws.Range("A103").Formula = '="01-04-2021"'
$ws.Range("A103").Copy()
$ws.Range("A103").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# --- Revised values for the existing "01-01-2021" row (row 102) ---
$row102 = @(54616, 41167, 33170, 3482, 14728, 14961, 7997, 10810, 6493, 4317, 2639, 17252, 16099, 2103, 9495, 8757, 738, 4501, 1153, 15434, 14025, 310, 1234, 12481, 1409, 56433)

# --- Values for the new "01-04-2021" row (row 103) ---
$row103 = @(55836, 43315, 34821, 4356, 14720, 15745, 8494, 11578, 7135, 4443, 943, 17716, 16617, 1001, 10658, 9755, 903, 4958, 1099, 16126, 14766, 318, 1384, 13063, 1360, 57427)

for ($i = 0; $i -lt $row102.Length; $i++) {
    $col = $i + 2  # Column B is index 2
    $ws.Cells.Item(102, $col).Value2 = $row102[$i]
    $ws.Cells.Item(103, $col).Value2 = $row103[$i]
}
